# Update the "want to go" counts (column F) on both the "展览" sheet and
# the "全部类型" sheet, which carry duplicate rows of event data.
$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 350
    $ws.Range("F4").Value = 74
}
